$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("List1")
$ws2 = $wb.Worksheets.Item("Měření aktivity")

# ---------------------------------------------------------------------------
# List1: correct the A51 measurement timestamp (time component fixed from
# 14:45 to 14:27). Every formula that depends on it (B51, C51, C52, E50)
# recalculates automatically.
# ---------------------------------------------------------------------------
$ws1.Range("A51").Value = 45530.602083333331

# ---------------------------------------------------------------------------
# "Měření aktivity": fill in the two measurement rows (27 & 28) that were
# previously left blank placeholders in the table, then grow the table by
# two more blank rows (30 & 31) the same way Excel does automatically when
# new data is typed right under a table.
# ---------------------------------------------------------------------------

# Row 27 data
$ws2.Range("A27").Value = 45526.636111111111
$ws2.Range("B27").Value = 0.006
$ws2.Range("C27").Value = 0.004
$ws2.Range("D27").Value = 0.004
$ws2.Range("E27").Value = 0.004
$ws2.Range("F27").Value = 0.004
$ws2.Range("G27").Value = 14.27
$ws2.Range("H27").Value = 14.28
$ws2.Range("I27").Value = 14.27
$ws2.Range("J27").Value = 14.27
$ws2.Range("K27").Value = 14.28
$ws2.Range("L27").Value = 14.27
$ws2.Range("M27").Value = 14.28
$ws2.Range("N27").Value = 14.28
$ws2.Range("O27").Value = 14.27
$ws2.Range("P27").Value = 14.27

# Row 28 data
$ws2.Range("A28").Value = 45530.602083333331
$ws2.Range("B28").Value = 0.009
$ws2.Range("C28").Value = 0.007
$ws2.Range("D28").Value = 0.006
$ws2.Range("E28").Value = 0.006
$ws2.Range("F28").Value = 0.006
$ws2.Range("G28").Value = 10.13
$ws2.Range("H28").Value = 10.13
$ws2.Range("I28").Value = 10.12
$ws2.Range("J28").Value = 10.13
$ws2.Range("K28").Value = 10.13
$ws2.Range("L28").Value = 10.13
$ws2.Range("M28").Value = 10.13
$ws2.Range("N28").Value = 10.13
$ws2.Range("O28").Value = 10.13
$ws2.Range("P28").Value = 10.13

# Grow the table ("Tabulka1") by two rows (30 & 31), mirroring the borders
# that previously sat on rows 27/28 (plain inner row) and row 29 (the
# bottom/closing row of the table) respectively, then resize the table
# definition + AutoFilter to cover the new range.
$ws2.Range("A27:Q27").Copy()
$ws2.Range("A30:Q30").PasteSpecial(-4122)
$ws2.Range("A29:Q29").Copy()
$ws2.Range("A31:Q31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("Q30").Formula = "=AVERAGE(G30:P30) - AVERAGE(Tabulka1[[#This Row],[č.1]:[č.5]])"
$ws2.Range("Q31").Formula = "=AVERAGE(G31:P31) - AVERAGE(Tabulka1[[#This Row],[č.1]:[č.5]])"

$lo = $ws2.ListObjects.Item("Tabulka1")
$lo.Resize($ws2.Range("A6:Q31"))

# ---------------------------------------------------------------------------
# Restore the selection state recorded in the saved view for each sheet.
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("A52").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("Q28").Select() | Out-Null
